# Actualización automática del mapa PEBCOM
# El caso "-560" (Pinzon 1590) de la fila 67 fue resuelto/eliminado del
# reporte. Se elimina esa fila completa y el resto de los casos (filas
# 68:85) se recorren una posición hacia arriba, tal como lo hace Excel
# al borrar una fila entera.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(67).Delete()
